$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $val) {
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $val
    $cellRange.Style = "Normal"
}

$ws.Range("D2").Value = '61.219.62'
$ws.Range("E2").Value = '  -0.24%  '

$ws.Range("D3").Value = '3.374.35'
$ws.Range("E3").Value = '  +1.25%  '

Set-TextValue $ws.Range("D4") '0.999'
$ws.Range("E4").Value = '  -0.06%  '

Set-TextValue $ws.Range("D5") '572.23'
$ws.Range("E5").Value = '  +0.47%  '

Set-TextValue $ws.Range("D6") '137.59'
$ws.Range("E6").Value = '  +7.20%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").Value = '3.373.85'
$ws.Range("E8").Value = '  +1.20%  '

Set-TextValue $ws.Range("D9") '0.477'
$ws.Range("E9").Value = '  -0.15%  '

Set-TextValue $ws.Range("D10") '7.60'
$ws.Range("E10").Value = '  +2.83%  '

Set-TextValue $ws.Range("D11") '0.124'
$ws.Range("E11").Value = '  +4.62%  '

Set-TextValue $ws.Range("D12") '0.394'
$ws.Range("E12").Value = '  +4.96%  '

$ws.Range("D13").Value = '3.949.88'
$ws.Range("E13").Value = '  +1.45%  '

$ws.Range("E14").Value = '  +1.78%  '

Set-TextValue $ws.Range("D15") '0.0000174'
$ws.Range("E15").Value = '  +3.05%  '

$ws.Range("D16").Value = '3.369.04'
$ws.Range("E16").Value = '  +1.40%  '

Set-TextValue $ws.Range("D17") '25.20'
$ws.Range("E17").Value = '  +1.15%  '

$ws.Range("D18").Value = '61.271.62'
$ws.Range("E18").Value = '  -0.28%  '

$ws.Range("E19").Value = '  +4.72%  '

Set-TextValue $ws.Range("D20") '5.87'
$ws.Range("E20").Value = '  +4.13%  '

Set-TextValue $ws.Range("D21") '9.39'
$ws.Range("E21").Value = '  +3.89%  '

Set-TextValue $ws.Range("D22") '379.34'
$ws.Range("E22").Value = '  +6.17%  '

Set-TextValue $ws.Range("D23") '0.569'
$ws.Range("E23").Value = '  +2.62%  '

$ws.Range("D24").Value = '3.508.30'
$ws.Range("E24").Value = '  +1.40%  '

$ws.Range("E25").Value = '  +0.03%  '

Set-TextValue $ws.Range("D26") '70.71'
$ws.Range("E26").Value = '  +1.24%  '

Set-TextValue $ws.Range("D27") '0.0000120'
$ws.Range("E27").Value = '  +11.07%  '

$ws.Range("E28").Value = '  +15.69%  '

Set-TextValue $ws.Range("D29") '7.77'
$ws.Range("E29").Value = '  +7.53%  '

$ws.Range("E30").Value = '  +0.23%  '

Set-TextValue $ws.Range("D31") '8.20'
$ws.Range("E31").Value = '  +3.39%  '

$ws.Range("E32").Value = '  +4.39%  '

Set-TextValue $ws.Range("D33") '2.12'
$ws.Range("E33").Value = '  +0.37%  '

$ws.Range("E34").Value = '  -0.02%  '

$ws.Range("D35").Value = '3.403.45'
$ws.Range("E35").Value = '  +1.34%  '

Set-TextValue $ws.Range("D36") '23.39'
$ws.Range("E36").Value = '  +3.62%  '

Set-TextValue $ws.Range("D37") '5.61'
$ws.Range("E37").Value = '  +5.69%  '

Set-TextValue $ws.Range("D38") '7.05'
$ws.Range("E38").Value = '  +4.22%  '

$ws.Range("E39").Value = '  +4.10%  '

Set-TextValue $ws.Range("D40") '162.65'
$ws.Range("E40").Value = '  +0.65%  '

Set-TextValue $ws.Range("D41") '0.0797'
$ws.Range("E41").Value = '  +5.07%  '

Set-TextValue $ws.Range("D42") '0.999'
$ws.Range("E42").Value = '  -0.14%  '

$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range("D43") '1.72'
$ws.Range("E43").Value = '  +8.85%  '

$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range("D44") '4.44'
$ws.Range("E44").Value = '  +1.81%  '

$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range("D45") '41.49'
$ws.Range("E45").Value = '  +0.40%  '

$ws.Range("B46").Value = 'ONDO'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextValue $ws.Range("D46") '1.20'
$ws.Range("E46").Value = '  +6.45%  '

$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range("D47") '0.759'
$ws.Range("E47").Value = '  +1.55%  '

Set-TextValue $ws.Range("D48") '23.28'
$ws.Range("E48").Value = '  +4.47%  '

Set-TextValue $ws.Range("D49") '6.96'
$ws.Range("E49").Value = '  +3.95%  '

Set-TextValue $ws.Range("D50") '22.97'
$ws.Range("E50").Value = '  +8.07%  '

$ws.Range("D51").Value = '2.327.82'
$ws.Range("E51").Value = '  +6.30%  '
